# Update "想去人数" (want-to-go count) figures in the "展览" (Exhibition)
# and "全部类型" (All types) sheets, plus the sold-out status text for row 25.

$wb = $excel.ActiveWorkbook

# Values shared by both the "展览" and "全部类型" sheets (column F, numeric).
$commonUpdates = @{
    2  = 142
    3  = 54
    5  = 98
    7  = 1307
    8  = 1554
    10 = 420
    12 = 171
    14 = 70
    15 = 112
    16 = 276
    17 = 315
    19 = 1755
    22 = 178
    23 = 684
    25 = 343
    26 = 4225
    28 = 285
    29 = 1109
    30 = 494
    36 = 147
}

# Sheet-specific values (these two rows differ by 1 between the sheets).
$sheetSpecificUpdates = @{
    "展览"   = @{ 32 = 597; 34 = 288 }
    "全部类型" = @{ 32 = 598; 34 = 289 }
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value2 = $commonUpdates[$row]
    }

    $specific = $sheetSpecificUpdates[$sheetName]
    foreach ($row in $specific.Keys) {
        $ws.Cells.Item($row, 6).Value2 = $specific[$row]
    }

    # Row 25's status text changed from "已售罄" to "暂时售罄".
    $ws.Cells.Item(25, 7).Value = "暂时售罄"
}
